$wb = $excel.ActiveWorkbook

$newGuid = "d7caa322-5ac1-430e-ba12-0ef535f71a30"
$newFileName = "$newGuid.md"
$newPath = "e2e\$newGuid.md"
$newHoDate = "2016-08-18 16:56:23"
$newZhFile = "$newGuid.5a7239e98103a6ce42c8d111a00091670c71668c.zh-cn.xlf"
$newZhDate = "2016-08-18 16:56:17"
$newDeFile = "$newGuid.5a7239e98103a6ce42c8d111a00091670c71668c.de-de.xlf"
$newDeDate = "2016-08-18 16:56:23"

$oldUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b2bf1ab1381840269cd18d74efbc92a68e76f9f/e2e/c674612d-3dbc-4443-944d-f89628fe0b6d.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("G2").Value = $newHoDate

$rngB2 = $wsOverview.Range("B2")
$rngB2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($rngB2, $oldUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newPath)

# --- zh-cn sheet ---
$rngZhA2 = $wsZhCn.Range("A2")
$rngZhA2.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($rngZhA2, $oldUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newFileName)
$wsZhCn.Range("G2").Value = $newZhFile
$wsZhCn.Range("H2").Value = $newZhDate

# --- de-de sheet ---
$rngDeA2 = $wsDeDe.Range("A2")
$rngDeA2.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($rngDeA2, $oldUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newFileName)
$wsDeDe.Range("G2").Value = $newDeFile
$wsDeDe.Range("H2").Value = $newHoDate
